$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the top of the data block (old rows
# 57-59 shift down to 58-60); insert a fresh row 57 and fill it in.
$ws.Rows(57).Insert()

$ws.Cells.Item(57, 1).Value = 4
$ws.Cells.Item(57, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(57, 3).Value = "Los Lagos"
$ws.Cells.Item(57, 4).Value = 44610
$ws.Cells.Item(57, 5).Value = 10
$ws.Cells.Item(57, 6).Value = 100112031
$ws.Cells.Item(57, 7).Value = "Poroto verde"
$ws.Cells.Item(57, 8).Value = "Magnum"
$ws.Cells.Item(57, 9).Value = "Primera"
$ws.Cells.Item(57, 10).Value = 50
$ws.Cells.Item(57, 11).Value = 30000
$ws.Cells.Item(57, 12).Value = 30000
$ws.Cells.Item(57, 13).Value = 30000
$ws.Cells.Item(57, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(57, 15).Value = "Región Metropolitana"
$ws.Cells.Item(57, 16).Value = 1200
$ws.Cells.Item(57, 17).Value = 25
$ws.Cells.Item(57, 18).Value = "Hortaliza"
